$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-26, per the commit's regenerated
# save_data (K replaces Strike#, std/mean recalculated, s_vals rewritten).
$newValues = @{
    2  = 0
    3  = 0
    4  = 3
    5  = 0
    6  = 1
    7  = 0
    8  = 2
    9  = 3
    10 = 5
    11 = 4
    12 = 3
    13 = 7
    14 = 1
    15 = 2
    16 = 2
    17 = 5
    18 = 4
    19 = 4
    20 = 6
    21 = 5
    22 = 5
    23 = 5
    24 = 6
    25 = 0
    26 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
